$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A48").Value = "What the maximum number of data files I can load?"
$ws.Range("B48").Value = "The maximum number of data files you can load is unlimited."
